$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 40.285713
$ws.Range("I6").Value = 40.285713
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 120.857139
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -8.857139000000004
$ws.Range("N6").Value = $null

# row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 112.5
$ws.Range("I8").Value = 112.5
$ws.Range("K8").Value = 337.5
$ws.Range("M8").Value = -198.5

# row 57 (Leve Item ID 43247)
$ws.Range("H57").Value = 40000
$ws.Range("J57").Value = 40000
$ws.Range("L57").Value = 120000
$ws.Range("N57").Value = -120998

# row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 1750
$ws.Range("J70").Value = 1750
$ws.Range("L70").Value = 5250
$ws.Range("N70").Value = -5790

# row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 1750
$ws.Range("J73").Value = 1750
$ws.Range("L73").Value = 5250
$ws.Range("N73").Value = -7122

# row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 1633.1666
$ws.Range("I80").Value = 1200.1
$ws.Range("K80").Value = 3600.3
$ws.Range("M80").Value = -2602.3

# row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 1633.1666
$ws.Range("I83").Value = 1200.1
$ws.Range("K83").Value = 10800.9
$ws.Range("M83").Value = -5808.9

# row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 12444.333
$ws.Range("J116").Value = 3666.5
$ws.Range("L116").Value = 3666.5
$ws.Range("N116").Value = -10550.5

# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 7245.625
$ws.Range("I132").Value = 7723.5713
$ws.Range("K132").Value = 23170.7139
$ws.Range("M132").Value = -20640.7139

# row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 1643.9166
$ws.Range("I135").Value = 1709.2
$ws.Range("K135").Value = 15382.8
$ws.Range("M135").Value = -12847.8

# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

# row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 999.6667
$ws.Range("I141").Value = 999.6667
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2999.0001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2180.9999
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1112.3334
$ws.Range("I2").Value = 1112.3334
$ws.Range("K2").Value = 1112.3334
$ws.Range("M2").Value = -999.3334

# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 1851.2
$ws.Range("I32").Value = 1328.3077
$ws.Range("K32").Value = 1328.3077
$ws.Range("M32").Value = -1041.3077

# row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3749.5
$ws.Range("I45").Value = 3999
$ws.Range("K45").Value = 3999
$ws.Range("M45").Value = -3622

# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1112.3334
$ws.Range("I116").Value = 1112.3334
$ws.Range("K116").Value = 1112.3334
$ws.Range("M116").Value = 1181.6666

$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1112.3334
$ws.Range("I3").Value = 1112.3334
$ws.Range("K3").Value = 1112.3334
$ws.Range("M3").Value = -998.3334

# row 75 (Leve Item ID 11872)
$ws.Range("H75").Value = 43271.285
$ws.Range("I75").Value = 16579.8
$ws.Range("J75").Value = 110000
$ws.Range("K75").Value = 16579.8
$ws.Range("L75").Value = 110000
$ws.Range("M75").Value = -15643.8
$ws.Range("N75").Value = -111872

# row 78 (Leve Item ID 11872)
$ws.Range("H78").Value = 43271.285
$ws.Range("I78").Value = 16579.8
$ws.Range("J78").Value = 110000
$ws.Range("K78").Value = 49739.39999999999
$ws.Range("L78").Value = 330000
$ws.Range("M78").Value = -45059.39999999999
$ws.Range("N78").Value = -339360

$ws = $wb.Worksheets.Item("CRP")
# row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 4653.933
$ws.Range("I58").Value = 1726.25
$ws.Range("J58").Value = 7999.857
$ws.Range("K58").Value = 1726.25
$ws.Range("L58").Value = 7999.857
$ws.Range("M58").Value = -1523.25
$ws.Range("N58").Value = -8405.857

# row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 4653.933
$ws.Range("I136").Value = 1726.25
$ws.Range("J136").Value = 7999.857
$ws.Range("K136").Value = 5178.75
$ws.Range("L136").Value = 23999.571
$ws.Range("M136").Value = -2628.75
$ws.Range("N136").Value = -29099.571

$ws = $wb.Worksheets.Item("CUL")
# row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 869.3333
$ws.Range("I5").Value = 714
$ws.Range("J5").Value = 947
$ws.Range("K5").Value = 2142
$ws.Range("L5").Value = 2841
$ws.Range("M5").Value = -2030
$ws.Range("N5").Value = -3065

# row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null

# row 59 (Leve Item ID 4694)
$ws.Range("H59").Value = 1118.3334
$ws.Range("I59").Value = 677.5
$ws.Range("K59").Value = 2032.5
$ws.Range("M59").Value = -1492.5

# row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 869.3333
$ws.Range("I135").Value = 714
$ws.Range("J135").Value = 947
$ws.Range("K135").Value = 6426
$ws.Range("L135").Value = 8523
$ws.Range("M135").Value = -3891
$ws.Range("N135").Value = -13593

$ws = $wb.Worksheets.Item("GSM")
# row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 4999
$ws.Range("K70").Value = 4999
$ws.Range("M70").Value = -4729

# row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 4999
$ws.Range("K73").Value = 4999
$ws.Range("M73").Value = -4063

# row 95 (Leve Item ID 18235)
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

# row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1874.75
$ws.Range("J102").Value = 1833.3334
$ws.Range("L102").Value = 1833.3334
$ws.Range("N102").Value = -5077.3334

# row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 3798.5789
$ws.Range("I113").Value = 4439.143
$ws.Range("J113").Value = 2005
$ws.Range("K113").Value = 4439.143
$ws.Range("L113").Value = 2005
$ws.Range("M113").Value = -2269.143
$ws.Range("N113").Value = -6345

# row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1569.2
$ws.Range("I122").Value = 1569.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4707.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2257.6
$ws.Range("N122").Value = $null

# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 1165.9375
$ws.Range("I132").Value = 1165.9375
$ws.Range("K132").Value = 3497.8125
$ws.Range("M132").Value = -967.8125

$ws = $wb.Worksheets.Item("LTW")
# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2300
$ws.Range("I22").Value = 2300
$ws.Range("K22").Value = 2300
$ws.Range("M22").Value = -2005

# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2300
$ws.Range("I27").Value = 2300
$ws.Range("K27").Value = 2300
$ws.Range("M27").Value = -2193

# row 69 (Leve Item ID 10671)
$ws.Range("H69").Value = 60000
$ws.Range("J69").Value = 60000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61622

# row 72 (Leve Item ID 10671)
$ws.Range("H72").Value = 60000
$ws.Range("J72").Value = 60000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188112

# row 74 (Leve Item ID 11990)
$ws.Range("H74").Value = 90000
$ws.Range("I74").Value = 90000
$ws.Range("K74").Value = 90000
$ws.Range("M74").Value = -89002

# row 77 (Leve Item ID 11990)
$ws.Range("H77").Value = 90000
$ws.Range("I77").Value = 90000
$ws.Range("K77").Value = 270000
$ws.Range("M77").Value = -265008

# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null

# row 140 (Leve Item ID 42503)
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# row 69 (Leve Item ID 10951)
$ws.Range("H69").Value = 25237.25
$ws.Range("J69").Value = 25237.25
$ws.Range("L69").Value = 25237.25
$ws.Range("N69").Value = -26735.25

# row 72 (Leve Item ID 10951)
$ws.Range("H72").Value = 25237.25
$ws.Range("J72").Value = 25237.25
$ws.Range("L72").Value = 75711.75
$ws.Range("N72").Value = -83199.75
